$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers (Area / Atotal), plus repeated Atotal/Qtotal in H/J/K ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- B2 / C2 used to hold the text "-"; now they are numeric zeros ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- New "Area" column (G) and totals (H, J, K) ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Match the new active selection (D2) recorded in the worksheet view ---
$ws.Range("D2").Select() | Out-Null
